# "merge all working test cases over"
# Adds a new "NewUrl" column (N) to the working test-case sheet:
#   - N1 header = "NewUrl"
#   - N2 value  = "Yes"
# Also backfills a previously-empty BundleQty cell (H2) with 2, and
# updates the view state (scroll/selection) to match the edited sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column N: header + the one data row.
$ws.Range("N1").Value = "NewUrl"
$ws.Range("N2").Value = "Yes"

# BundleQty (column H) was blank on row 2; fill in the quantity.
$ws.Cells.Item(2, 8).Value = 2

# Scroll the view over a column and leave the selection on the new
# column's second row, matching the saved view state of the edit.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("N3").Select()
